# Apply updated cryptocurrency price/volume figures to sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.200.79"
$ws.Range("E2").Value = "  -4.70%  "
$ws.Range("D3").Value = "2.236.03"
$ws.Range("E3").Value = "  -5.57%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'318.92"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "'100.81"
$ws.Range("E6").Value = "  -6.62%  "
$ws.Range("E7").Value = "  -7.21%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.563"
$ws.Range("E9").Value = "  -8.07%  "
$ws.Range("D10").Value = "'36.98"
$ws.Range("E10").Value = "  -9.30%  "
$ws.Range("D11").Value = "'54.46"
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("D12").Value = "'0.0826"
$ws.Range("E12").Value = "  -9.85%  "
$ws.Range("E13").Value = "  -9.28%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "2.575.78"
$ws.Range("E15").Value = "  -5.60%  "
$ws.Range("D16").Value = "'0.862"
$ws.Range("E16").Value = "  -12.16%  "
$ws.Range("D17").Value = "'14.36"
$ws.Range("E17").Value = "  -6.27%  "
$ws.Range("D18").Value = "2.237.63"
$ws.Range("E18").Value = "  -5.06%  "
$ws.Range("D19").Value = "43.127.82"
$ws.Range("E19").Value = "  -4.79%  "
$ws.Range("E20").Value = "  -6.67%  "
$ws.Range("D21").Value = "0.0₃0966"
$ws.Range("E21").Value = "  -9.03%  "
$ws.Range("E22").Value = "  -10.35%  "
$ws.Range("D23").Value = "'65.42"
$ws.Range("E23").Value = "  -10.74%  "
$ws.Range("D24").Value = "'3.18"
$ws.Range("E24").Value = "  -11.57%  "
$ws.Range("D25").Value = "'238.19"
$ws.Range("E25").Value = "  -8.77%  "
$ws.Range("E26").Value = "  -8.16%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'4.08"
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").Value = "'10.06"
$ws.Range("E29").Value = "  -9.80%  "
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").Value = "'6.39"
$ws.Range("E31").Value = "  -14.68%  "
$ws.Range("D32").Value = "'35.48"
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("D33").Value = "'20.49"
$ws.Range("E33").Value = "  -8.12%  "
$ws.Range("E34").Value = "  -9.30%  "
$ws.Range("D35").Value = "'153.42"
$ws.Range("D36").Value = "'2.75"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "'3.12"
$ws.Range("E37").Value = "  +6.66%  "
$ws.Range("D38").Value = "'1.97"
$ws.Range("E38").Value = "  +4.76%  "
$ws.Range("E39").Value = "  -6.87%  "
$ws.Range("E40").Value = "  -4.90%  "
$ws.Range("E41").Value = "  -11.02%  "
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = "  -7.18%  "
$ws.Range("E43").Value = "  -8.12%  "
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "1.799.64"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'87.11"
$ws.Range("E47").Value = "  -11.61%  "
$ws.Range("E48").Value = "  -9.62%  "
$ws.Range("D49").Value = "'76.50"
$ws.Range("E49").Value = "  -7.38%  "
$ws.Range("E50").Value = "  -10.20%  "
$ws.Range("D51").Value = "'59.31"
$ws.Range("E51").Value = "  -15.36%  "

# The apostrophe-prefix above forces Excel to keep these numeric-looking
# strings as text (matching the source data), but it also stamps a
# "quote prefix" style on the cell. Reset style back to Normal so the
# cells keep their original (unstyled) appearance.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
